$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header text updates (shared strings): volume number + report week dates
# ---------------------------------------------------------------------------
$ws.Range("A8").Value2 = "Volume 30   Number  35"
$ws.Range("C9").Value2 = "Report Covering the Week  8/28/2023  Through  9/3/2023"

# ---------------------------------------------------------------------------
# 2) Structural (type/style) fix-ups: cells that flip between a numeric value
#    and the text markers "0" / "***.*" (shared strings already used
#    elsewhere on the sheet). Copying an existing matching cell reproduces
#    both the text-type and the original cell style faithfully; the actual
#    display values are then applied on top where they differ from the
#    source cell that was copied.
# ---------------------------------------------------------------------------

# C14 : number 1  -> text "0"
$ws.Range("D14").Copy($ws.Range("C14"))

# C26 : number 1  -> text "0"
$ws.Range("D14").Copy($ws.Range("C26"))

# C27 : number 3  -> text "0"
$ws.Range("D14").Copy($ws.Range("C27"))

# D28 : number 1  -> text "0"
$ws.Range("D14").Copy($ws.Range("D28"))

# E28 : number -100 -> text "***.*"
$ws.Range("E14").Copy($ws.Range("E28"))

# D29 : number 1  -> text "0"
$ws.Range("D14").Copy($ws.Range("D29"))

# E29 : number -100 -> text "***.*"
$ws.Range("E14").Copy($ws.Range("E29"))

# D23 : text "0"     -> number 1
$ws.Range("C23").Copy($ws.Range("D23"))
$ws.Range("D23").Value2 = 1

# E23 : text "***.*" -> number 0
$ws.Range("H23").Copy($ws.Range("E23"))
$ws.Range("E23").Value2 = 0

# ---------------------------------------------------------------------------
# 3) Plain numeric value updates
# ---------------------------------------------------------------------------

# Row 14
$ws.Range("M14").Value2 = -76.923076923076
$ws.Range("N14").Value2 = -85

# Row 15
$ws.Range("F15").Value2 = 1
$ws.Range("H15").Value2 = 0
$ws.Range("M15").Value2 = 0
$ws.Range("N15").Value2 = -67.924528301886

# Row 16
$ws.Range("D16").Value2 = 5
$ws.Range("E16").Value2 = -80
$ws.Range("F16").Value2 = 5
$ws.Range("G16").Value2 = 16
$ws.Range("H16").Value2 = -68.75
$ws.Range("I16").Value2 = 116
$ws.Range("J16").Value2 = 109
$ws.Range("K16").Value2 = 6.422018348623
$ws.Range("L16").Value2 = 8.411214953271
$ws.Range("M16").Value2 = -47.511312217194
$ws.Range("N16").Value2 = -85.572139303482

# Row 17
$ws.Range("D17").Value2 = 1
$ws.Range("E17").Value2 = 100
$ws.Range("F17").Value2 = 20
$ws.Range("G17").Value2 = 20
$ws.Range("I17").Value2 = 216
$ws.Range("J17").Value2 = 228
$ws.Range("K17").Value2 = -5.263157894736
$ws.Range("L17").Value2 = -1.369863013698
$ws.Range("M17").Value2 = -2.262443438914
$ws.Range("N17").Value2 = -65.273311897106

# Row 18
$ws.Range("C18").Value2 = 4
$ws.Range("D18").Value2 = 3
$ws.Range("E18").Value2 = 33.333333333333
$ws.Range("F18").Value2 = 9
$ws.Range("H18").Value2 = -25
$ws.Range("I18").Value2 = 108
$ws.Range("J18").Value2 = 132
$ws.Range("K18").Value2 = -18.181818181818
$ws.Range("L18").Value2 = 2.857142857142
$ws.Range("M18").Value2 = -25
$ws.Range("N18").Value2 = -79.545454545454

# Row 19
$ws.Range("C19").Value2 = 3
$ws.Range("D19").Value2 = 11
$ws.Range("E19").Value2 = -72.727272727272
$ws.Range("F19").Value2 = 21
$ws.Range("H19").Value2 = -27.586206896551
$ws.Range("I19").Value2 = 226
$ws.Range("J19").Value2 = 275
$ws.Range("K19").Value2 = -17.818181818181
$ws.Range("L19").Value2 = -6.995884773662
$ws.Range("M19").Value2 = 8.653846153846
$ws.Range("N19").Value2 = -9.6

# Row 20
$ws.Range("C20").Value2 = 4
$ws.Range("D20").Value2 = 4
$ws.Range("E20").Value2 = 0
$ws.Range("F20").Value2 = 12
$ws.Range("G20").Value2 = 10
$ws.Range("H20").Value2 = 20
$ws.Range("I20").Value2 = 72
$ws.Range("J20").Value2 = 94
$ws.Range("K20").Value2 = -23.404255319148
$ws.Range("L20").Value2 = 9.090909090909
$ws.Range("M20").Value2 = 20
$ws.Range("N20").Value2 = -82.608695652173

# Row 21
$ws.Range("C21").Value2 = 14
$ws.Range("D21").Value2 = 24
$ws.Range("E21").Value2 = -41.666666666666
$ws.Range("F21").Value2 = 69
$ws.Range("G21").Value2 = 88
$ws.Range("H21").Value2 = -21.590909090909
$ws.Range("I21").Value2 = 758
$ws.Range("J21").Value2 = 858
$ws.Range("K21").Value2 = -11.655011655011
$ws.Range("L21").Value2 = 0.132100396301
$ws.Range("M21").Value2 = -14.253393665158
$ws.Range("N21").Value2 = -71.832032701597

# Row 22
$ws.Range("M22").Value2 = -44.444444444444

# Row 23 (D23/E23 type-fix handled above)
$ws.Range("F23").Value2 = 7
$ws.Range("G23").Value2 = 6
$ws.Range("H23").Value2 = 16.666666666666
$ws.Range("I23").Value2 = 60
$ws.Range("J23").Value2 = 55
$ws.Range("K23").Value2 = 9.090909090909
$ws.Range("L23").Value2 = -18.918918918918
$ws.Range("M23").Value2 = 5.263157894736

# Row 24
$ws.Range("C24").Value2 = 15
$ws.Range("D24").Value2 = 19
$ws.Range("E24").Value2 = -21.052631578947
$ws.Range("F24").Value2 = 75
$ws.Range("G24").Value2 = 90
$ws.Range("H24").Value2 = -16.666666666666
$ws.Range("I24").Value2 = 580
$ws.Range("J24").Value2 = 555
$ws.Range("K24").Value2 = 4.504504504504
$ws.Range("L24").Value2 = 50.649350649350
$ws.Range("M24").Value2 = 8.208955223880

# Row 25
$ws.Range("C25").Value2 = 13
$ws.Range("D25").Value2 = 10
$ws.Range("E25").Value2 = 30
$ws.Range("F25").Value2 = 31
$ws.Range("G25").Value2 = 35
$ws.Range("H25").Value2 = -11.428571428571
$ws.Range("I25").Value2 = 349
$ws.Range("J25").Value2 = 281
$ws.Range("K25").Value2 = 24.199288256227
$ws.Range("L25").Value2 = 53.744493392070
$ws.Range("M25").Value2 = -35.370370370370

# Row 26 (C26 type-fix handled above)
$ws.Range("F26").Value2 = 2
$ws.Range("H26").Value2 = 100

# Row 27 (C27 type-fix handled above) - no other numeric changes on row 27

# Row 28 (D28/E28 type-fix handled above)
$ws.Range("F28").Value2 = 1
$ws.Range("H28").Value2 = -66.666666666666
$ws.Range("M28").Value2 = -72.5
$ws.Range("N28").Value2 = -91.729323308270

# Row 29 (D29/E29 type-fix handled above)
$ws.Range("F29").Value2 = 1
$ws.Range("H29").Value2 = -66.666666666666
$ws.Range("M29").Value2 = -74.285714285714
$ws.Range("N29").Value2 = -92.372881355932
